$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------
# Sheet2: Question 1 - Square table (rows 1-7)
# ---------------------------------------------------------------
$ws2.Range("A1").Value = "Question 1"
$ws2.Range("A1").Font.Bold = $true

$ws2.Range("B2").Value = "Number"
$ws2.Range("C2").Value = "Square"
$ws2.Range("B2:C2").HorizontalAlignment = -4108

$ws2.Range("B3").Value = 8
$ws2.Range("B4").Value = 12
$ws2.Range("B5").Value = 14
$ws2.Range("B6").Value = 18
$ws2.Range("B7").Value = 24
$ws2.Range("B3:B7").HorizontalAlignment = -4108

$ws2.Range("C3").Formula = "=POWER(B3,2)"
$ws2.Range("C4:C7").Formula = "=POWER(B4,2)"
$ws2.Range("C3:C7").HorizontalAlignment = -4108

$r1 = $ws2.Range("B2:C7")
$r1.Borders.LineStyle = 1
$r1.Borders.Weight = 2
$r1.BorderAround(1, -4138)

# ---------------------------------------------------------------
# Sheet2: Question 2 - Cube table (rows 9-15)
# ---------------------------------------------------------------
$ws2.Range("A9").Value = "Question 2"

$ws2.Range("B10").Value = "Number"
$ws2.Range("C10").Value = "Cube"
$ws2.Range("B10:C10").HorizontalAlignment = -4108

$ws2.Range("B11").Value = 8
$ws2.Range("B12").Value = 12
$ws2.Range("B13").Value = 14
$ws2.Range("B14").Value = 18
$ws2.Range("B15").Value = 24
$ws2.Range("B11:B15").HorizontalAlignment = -4108

$ws2.Range("C11").Formula = "=POWER(B11,3)"
$ws2.Range("C12:C15").Formula = "=POWER(B12,3)"
$ws2.Range("C11:C15").HorizontalAlignment = -4108

$r2 = $ws2.Range("B10:C15")
$r2.Borders.LineStyle = 1
$r2.Borders.Weight = 2
$r2.BorderAround(1, -4138)

# ---------------------------------------------------------------
# Sheet2: Question 3 - Circle table (rows 17-23)
# ---------------------------------------------------------------
$ws2.Range("A17").Value = "Question 3"

$ws2.Range("B18").Value = "Number"
$ws2.Range("C18").Value = "Radius"
$ws2.Range("D18").Value = "Circumference"
$ws2.Range("B18:D18").HorizontalAlignment = -4108

$ws2.Range("B19").Value = 8
$ws2.Range("B20").Value = 12
$ws2.Range("B21").Value = 14
$ws2.Range("B22").Value = 18
$ws2.Range("B23").Value = 24
$ws2.Range("B19:B23").HorizontalAlignment = -4108

$ws2.Range("C19").Formula = "=3.14*POWER(B19,2)"
$ws2.Range("C20:C23").Formula = "=3.14*POWER(B20,2)"

$ws2.Range("D19").Formula = "=2*3.14*B19"
$ws2.Range("D20:D23").Formula = "=2*3.14*B20"

$r3 = $ws2.Range("B18:D23")
$r3.Borders.LineStyle = 1
$r3.Borders.Weight = 2
$r3.BorderAround(1, -4138)

# ---------------------------------------------------------------
# Sheet2: Question 4 - Employee salary table (rows 25-29)
# ---------------------------------------------------------------
$ws2.Range("A25").Value = "Question 4"

$ws2.Range("B26").Value = "Employee Code"
$ws2.Range("C26").Value = "Name"
$ws2.Range("D26").Value = "Basic Salary"
$ws2.Range("E26").Value = "HRA"
$ws2.Range("F26").Value = "DA"
$ws2.Range("G26").Value = "CCA"
$ws2.Range("H26").Value = "GS"
$ws2.Range("I26").Value = "NS"

$ws2.Range("B27").Value = "Wz-0001"
$ws2.Range("C27").Value = "Ishminder"
$ws2.Range("D27").Value = 20000

$ws2.Range("B28").Value = "Wz-0002"
$ws2.Range("C28").Value = "Haresh"
$ws2.Range("D28").Value = 30000

$ws2.Range("B29").Value = "Wz-0003"
$ws2.Range("C29").Value = "Satyam"
$ws2.Range("D29").Value = 40000
$ws2.Range("C29").HorizontalAlignment = -4131
$ws2.Range("C29").WrapText = $true

$ws2.Range("E27").Formula = "=D27*40/100"
$ws2.Range("F27").Formula = "=10*D27/100"
$ws2.Range("G27").Formula = "=5*D27/100"
$ws2.Range("H27").Formula = "=SUM(D27:G27)"

$ws2.Range("E28:E29").Formula = "=D28*40/100"
$ws2.Range("F28:F29").Formula = "=10*D28/100"
$ws2.Range("G28:G29").Formula = "=5*D28/100"
$ws2.Range("H28:H29").Formula = "=SUM(D28:G28)"

# ---------------------------------------------------------------
# Column widths / view settings for Sheet2
# ---------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 10.5703125
$ws2.Columns.Item(2).ColumnWidth = 15
$ws2.Columns.Item(4).ColumnWidth = 14.140625

$ws2.Activate()
$ws2.Range("J27").Select()
$excel.ActiveWindow.ScrollRow = 7

# ---------------------------------------------------------------
# Sheet1 bookkeeping: random cells get re-evaluated (volatile funcs),
# just touch nothing else there.
# ---------------------------------------------------------------
